$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2995.4666
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 3152.2856
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 3152.2856
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -3804.2856
$ws.Range("H132").Value = 1067204.5
$ws.Range("I132").Value = 2054.7273
$ws.Range("K132").Value = 6164.1819
$ws.Range("M132").Value = -3634.1819
$ws.Range("H137").Value = 3521199.5
$ws.Range("I137").Value = 11040331
$ws.Range("J137").Value = 1605.9574
$ws.Range("K137").Value = 33120993
$ws.Range("L137").Value = 4817.8722
$ws.Range("M137").Value = -33118443
$ws.Range("N137").Value = -9917.8722

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1777.8889
$ws.Range("I2").Value = 1777.8889
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1777.8889
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1664.8889
$ws.Range("N2").ClearContents()
$ws.Range("H74").Value = 7876105.5
$ws.Range("I74").Value = 13945053
$ws.Range("J74").Value = 73173.28999999999
$ws.Range("K74").Value = 13945053
$ws.Range("L74").Value = 73173.28999999999
$ws.Range("M74").Value = -13944179
$ws.Range("N74").Value = -74921.28999999999
$ws.Range("H77").Value = 7876105.5
$ws.Range("I77").Value = 13945053
$ws.Range("J77").Value = 73173.28999999999
$ws.Range("K77").Value = 69725265
$ws.Range("L77").Value = 365866.45
$ws.Range("M77").Value = -69720897
$ws.Range("N77").Value = -374602.45
$ws.Range("H116").Value = 1777.8889
$ws.Range("I116").Value = 1777.8889
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1777.8889
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 516.1111000000001
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 11113382
$ws.Range("I122").Value = 2588.375
$ws.Range("J122").Value = 55556556
$ws.Range("K122").Value = 7765.125
$ws.Range("L122").Value = 166669668
$ws.Range("M122").Value = -5315.125
$ws.Range("N122").Value = -166674568
$ws.Range("H132").Value = 102584
$ws.Range("I132").Value = 102376
$ws.Range("J132").Value = 102792
$ws.Range("K132").Value = 307128
$ws.Range("L132").Value = 308376
$ws.Range("M132").Value = -304598
$ws.Range("N132").Value = -313436

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1777.8889
$ws.Range("I3").Value = 1777.8889
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1777.8889
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1663.8889
$ws.Range("N3").ClearContents()
$ws.Range("H76").Value = 38024.875
$ws.Range("J76").Value = 38024.875
$ws.Range("L76").Value = 38024.875
$ws.Range("N76").Value = -38654.875
$ws.Range("H79").Value = 38024.875
$ws.Range("J79").Value = 38024.875
$ws.Range("L79").Value = 38024.875
$ws.Range("N79").Value = -40208.875
$ws.Range("H94").Value = 603.03845
$ws.Range("I94").Value = 603.03845
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 603.03845
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -152.03845
$ws.Range("N94").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 38982.25
$ws.Range("J88").Value = 38982.25
$ws.Range("L88").Value = 38982.25
$ws.Range("N88").Value = -39794.25
$ws.Range("H91").Value = 38982.25
$ws.Range("J91").Value = 38982.25
$ws.Range("L91").Value = 38982.25
$ws.Range("N91").Value = -41790.25
$ws.Range("H99").Value = 4243.8335
$ws.Range("I99").Value = 4741.2
$ws.Range("J99").Value = 1757
$ws.Range("K99").Value = 4741.2
$ws.Range("L99").Value = 1757
$ws.Range("M99").Value = -3243.2
$ws.Range("N99").Value = -4753
$ws.Range("H122").Value = 1529.875
$ws.Range("J122").Value = 1172.25
$ws.Range("L122").Value = 3516.75
$ws.Range("N122").Value = -8416.75
$ws.Range("H126").Value = 4243.8335
$ws.Range("I126").Value = 4741.2
$ws.Range("J126").Value = 1757
$ws.Range("K126").Value = 14223.6
$ws.Range("L126").Value = 5271
$ws.Range("M126").Value = -11753.6
$ws.Range("N126").Value = -10211
$ws.Range("H141").Value = 130326
$ws.Range("J141").Value = 130326
$ws.Range("L141").Value = 130326
$ws.Range("N141").Value = -140686

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1200
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3600
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3938
$ws.Range("H23").Value = 87.59999999999999
$ws.Range("I23").Value = 46.6
$ws.Range("J23").Value = 128.6
$ws.Range("K23").Value = 139.8
$ws.Range("L23").Value = 385.8
$ws.Range("M23").Value = 95.19999999999999
$ws.Range("N23").Value = -855.8
$ws.Range("H34").Value = 758.1818
$ws.Range("I34").Value = 435
$ws.Range("J34").Value = 942.8570999999999
$ws.Range("K34").Value = 1305
$ws.Range("L34").Value = 2828.5713
$ws.Range("M34").Value = -1221
$ws.Range("N34").Value = -2996.5713
$ws.Range("H39").Value = 2925
$ws.Range("J39").Value = 2925
$ws.Range("L39").Value = 8775
$ws.Range("N39").Value = -9363
$ws.Range("H55").Value = 500
$ws.Range("J55").Value = 500
$ws.Range("L55").Value = 1500
$ws.Range("N55").Value = -1854
$ws.Range("H92").Value = 758.8889
$ws.Range("I92").Value = 750
$ws.Range("J92").Value = 770
$ws.Range("K92").Value = 2250
$ws.Range("L92").Value = 2310
$ws.Range("M92").Value = -1002
$ws.Range("N92").Value = -4806
$ws.Range("H121").Value = 204899860
$ws.Range("I121").Value = 1075
$ws.Range("J121").Value = 263442370
$ws.Range("K121").Value = 3225
$ws.Range("L121").Value = 790327110
$ws.Range("M121").Value = -1915
$ws.Range("N121").Value = -790329730

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1802.9412
$ws.Range("I113").Value = 1600.1428
$ws.Range("J113").Value = 1944.9
$ws.Range("K113").Value = 1600.1428
$ws.Range("L113").Value = 1944.9
$ws.Range("M113").Value = 569.8571999999999
$ws.Range("N113").Value = -6284.9
$ws.Range("H132").Value = 59836.145
$ws.Range("I132").Value = 52591.95
$ws.Range("J132").Value = 69495.07000000001
$ws.Range("K132").Value = 157775.85
$ws.Range("L132").Value = 208485.21
$ws.Range("M132").Value = -155245.85
$ws.Range("N132").Value = -213545.21

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 923.8421
$ws.Range("I46").Value = 604.0833
$ws.Range("J46").Value = 1472
$ws.Range("K46").Value = 604.0833
$ws.Range("L46").Value = 1472
$ws.Range("M46").Value = -416.0833
$ws.Range("N46").Value = -1848
$ws.Range("H61").Value = 3894.1304
$ws.Range("I61").Value = 3622.2856
$ws.Range("J61").Value = 4317
$ws.Range("K61").Value = 3622.2856
$ws.Range("L61").Value = 4317
$ws.Range("M61").Value = -3420.2856
$ws.Range("N61").Value = -4721
$ws.Range("H113").Value = 3894.1304
$ws.Range("I113").Value = 3622.2856
$ws.Range("J113").Value = 4317
$ws.Range("K113").Value = 3622.2856
$ws.Range("L113").Value = 4317
$ws.Range("M113").Value = -1452.2856
$ws.Range("N113").Value = -8657
$ws.Range("H122").Value = 3649.4707
$ws.Range("I122").Value = 3609.8572
$ws.Range("J122").Value = 3834.3333
$ws.Range("K122").Value = 10829.5716
$ws.Range("L122").Value = 11502.9999
$ws.Range("M122").Value = -8379.571599999999
$ws.Range("N122").Value = -16402.9999
$ws.Range("H132").Value = 41288.04
$ws.Range("I132").Value = 2454.8
$ws.Range("J132").Value = 65558.81
$ws.Range("K132").Value = 7364.400000000001
$ws.Range("L132").Value = 196676.43
$ws.Range("M132").Value = -4834.400000000001
$ws.Range("N132").Value = -201736.43

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 30000
$ws.Range("J16").Value = 30000
$ws.Range("L16").Value = 30000
$ws.Range("N16").Value = -30584
$ws.Range("H45").Value = 11400
$ws.Range("J45").Value = 11400
$ws.Range("L45").Value = 11400
$ws.Range("N45").Value = -12382
$ws.Range("H107").Value = 1258.25
$ws.Range("J107").Value = 1126.6666
$ws.Range("L107").Value = 3379.9998
$ws.Range("N107").Value = -7219.9998
$ws.Range("H108").Value = 37813.5
$ws.Range("J108").Value = 37813.5
$ws.Range("L108").Value = 37813.5
$ws.Range("N108").Value = -45493.5
$ws.Range("H113").Value = 1224.2916
$ws.Range("I113").Value = 1183.8422
$ws.Range("J113").Value = 1378
$ws.Range("K113").Value = 3551.5266
$ws.Range("L113").Value = 4134
$ws.Range("M113").Value = -1381.5266
$ws.Range("N113").Value = -8474
$ws.Range("H136").Value = 51657.227
$ws.Range("I136").Value = 36516.93
$ws.Range("J136").Value = 86984.586
$ws.Range("K136").Value = 109550.79
$ws.Range("L136").Value = 260953.758
$ws.Range("M136").Value = -107000.79
$ws.Range("N136").Value = -266053.758
